# Update "Home" row (row 2) stats on both the OFF and DEF sheets to reflect
# the additional divisional-round game logged for the simulated 2021 season.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 544
$wsOff.Range("C2").Value = 380
$wsOff.Range("D2").Value = 141
$wsOff.Range("E2").Value = 66
$wsOff.Range("G2").Value = 8

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 631
$wsDef.Range("C2").Value = 442
$wsDef.Range("D2").Value = 135
$wsDef.Range("E2").Value = 63

$wb.Save()
